# Test_Data.xlsx update:
# - Clear the "Add_Required" (F2) and "Is_Used" (G2) values from the second
#   registration row - those two columns/values are no longer used.
# - Duplicate the remaining registration row (A:E) into a new row 3, including
#   its mailto hyperlink on the Email cell, giving two identical test-data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the now-unused Add_Required / Is_Used values from row 2
$ws.Range("F2:G2").ClearContents()

# Add a second data row identical to row 2 (Name..Language)
$ws.Range("A2:E2").Copy($ws.Range("A3"))

# Row 3's Email cell (C3) needs its own hyperlink, same target as C2's
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:kapil@gmail.com")

# Copying the hyperlink style onto C3 pulled in Excel's default Hyperlink
# formatting; restore the same style C2 uses so both Email cells match
$ws.Range("C3").Style = $ws.Range("C2").Style

# Move the active selection the way the authored workbook left it
$ws.Range("E8").Select() | Out-Null
